# Trade #94 closed at 2026-02-17 15:54:50 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.81   # Current Capital
$wsSummary.Range("B4").Value = -0.2      # Total P&L $
$wsSummary.Range("B5").Value = -0.04     # Total P&L %
$wsSummary.Range("B6").Value = 94        # Total Trades
$wsSummary.Range("B7").Value = 34        # Winning Trades
$wsSummary.Range("B9").Value = 36.17     # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.81      # Capital
$wsStatus.Range("D4").Value = 94         # Trades
$wsStatus.Range("E4").Value = -0.2       # P&L $
$wsStatus.Range("F4").Value = -0.19      # P&L %
$wsStatus.Range("G4").Value = 36.17      # Win Rate %

# --- Add new trade row (#94) to both "All Trades" and "MarketMaking" sheets ---
$newRow = 95

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 94

    $ws.Range("B" + $row).NumberFormat = "@"
    $ws.Range("B" + $row).Value = "2026-02-17"

    $ws.Range("C" + $row).NumberFormat = "@"
    $ws.Range("C" + $row).Value = "15:54:43"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.9
    $ws.Cells.Item($row, 7).Value = 0.92
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 2.2222
    $ws.Cells.Item($row, 10).Value = 0.02
    $ws.Cells.Item($row, 11).Value = 99.81
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.11
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades $newRow

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking $newRow
